# Update "想去人数" (want-to-go count) in column F for the rows whose
# upstream source values changed, on both the "展览" and "全部类型" sheets
# (these two sheets mirror the same underlying data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2308
    3 = 1761
    6 = 945
    8 = 5872
    9 = 91
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
